$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 131 was missing its "End Time" entry - fill it in. The Duration /
# Second Duration / Absolute Value formulas already in D131:F131 recalc
# automatically once C131 has a value.
$ws.Range("C131").Value = 0.3756944444444445

# Row 139 already carried the Duration/Second Duration/Absolute Value
# formulas but was otherwise blank - fill in the Date/Start Time/End Time.
$ws.Range("A139").Value = 43461
$ws.Range("B139").Value = 0
$ws.Range("C139").Value = 0

# Insert two brand-new rows (140, 141) for the new daily records. Using
# "format from above" on the insert makes the new D/E/F cells pick up the
# same cell style used by row 139 right above them (matches the table's
# calculated-column formatting) before any content is written.
$ws.Range("A140:A141").EntireRow.Insert(-4121, 0)

# Row 140: new daily power record
$ws.Range("A140").Value = 43462
$ws.Range("B140").Value = 0
$ws.Range("C140").Value = 0
$ws.Range("D140").Formula = "=(C140-B140)* 1440"
$ws.Range("E140").Formula = "=IF(C140>B140, (C140-B140)*1440, (B140-C140)*1440)"
$ws.Range("F140").Formula = "=ABS((C140-B140)*1440)"

# Row 141: new daily power record
$ws.Range("A141").Value = 43463
$ws.Range("B141").Value = 0.33333333333333331
$ws.Range("C141").Value = 0.34722222222222227
$ws.Range("D141").Formula = "=(C141-B141)* 1440"
$ws.Range("E141").Formula = "=IF(C141>B141, (C141-B141)*1440, (B141-C141)*1440)"
$ws.Range("F141").Formula = "=ABS((C141-B141)*1440)"

# Extend the "comforter_cda_table" Excel table so it covers the two new rows.
$tbl = $ws.ListObjects.Item(1)
$tbl.Resize($ws.Range("A1:F141"))

# Leave the selection where the author ended up after the edit.
$ws.Range("D141").Select()
